$wb = $excel.ActiveWorkbook

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 14286434
$ws.Range("I103").Value = 443.33334
$ws.Range("J103").Value = 18182614
$ws.Range("K103").Value = 1330.00002
$ws.Range("L103").Value = 54547842
$ws.Range("M103").Value = -744.0000199999999
$ws.Range("N103").Value = -54549014

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4426.4
$ws.Range("I116").Value = 3289.1538
$ws.Range("J116").Value = 5296.0586
$ws.Range("K116").Value = 3289.1538
$ws.Range("L116").Value = 5296.0586
$ws.Range("M116").Value = 152.8462
$ws.Range("N116").Value = -12180.0586

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1794.262
$ws.Range("I137").Value = 891.86664
$ws.Range("J137").Value = 4050.25
$ws.Range("K137").Value = 2675.59992
$ws.Range("L137").Value = 12150.75
$ws.Range("M137").Value = -125.5999199999997
$ws.Range("N137").Value = -17250.75

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13176.76
$ws.Range("I32").Value = 4372.1157
$ws.Range("J32").Value = 32774.195
$ws.Range("K32").Value = 4372.1157
$ws.Range("L32").Value = 32774.195
$ws.Range("M32").Value = -4085.1157
$ws.Range("N32").Value = -33348.195

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1367.125
$ws.Range("I74").Value = 689.16327
$ws.Range("J74").Value = 6112.857
$ws.Range("K74").Value = 689.16327
$ws.Range("L74").Value = 6112.857
$ws.Range("M74").Value = 184.83673
$ws.Range("N74").Value = -7860.857

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1367.125
$ws.Range("I77").Value = 689.16327
$ws.Range("J77").Value = 6112.857
$ws.Range("K77").Value = 3445.81635
$ws.Range("L77").Value = 30564.285
$ws.Range("M77").Value = 922.1836499999999
$ws.Range("N77").Value = -39300.285

# ARM row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 42444
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 42444
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 42444
$ws.Range("N80").Value = -44440

# ARM row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 42444
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 42444
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 127332
$ws.Range("N83").Value = -137316

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 28520.375
$ws.Range("I82").Value = 9176.5
$ws.Range("J82").Value = 34968.332
$ws.Range("K82").Value = 9176.5
$ws.Range("L82").Value = 34968.332
$ws.Range("M82").Value = -8793.5
$ws.Range("N82").Value = -35734.332

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 28520.375
$ws.Range("I85").Value = 9176.5
$ws.Range("J85").Value = 34968.332
$ws.Range("K85").Value = 9176.5
$ws.Range("L85").Value = 34968.332
$ws.Range("M85").Value = -7850.5
$ws.Range("N85").Value = -37620.332

# CRP row 4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 100001.336
$ws.Range("I4").Value = 100000
$ws.Range("J4").Value = 100002
$ws.Range("K4").Value = 100000
$ws.Range("L4").Value = 100002
$ws.Range("M4").Value = -99888
$ws.Range("N4").Value = -100226

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3211.238
$ws.Range("I31").Value = 1474.4
$ws.Range("J31").Value = 4790.1816
$ws.Range("K31").Value = 1474.4
$ws.Range("L31").Value = 4790.1816
$ws.Range("M31").Value = -1179.4
$ws.Range("N31").Value = -5380.1816

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3211.238
$ws.Range("I34").Value = 1474.4
$ws.Range("J34").Value = 4790.1816
$ws.Range("K34").Value = 1474.4
$ws.Range("L34").Value = 4790.1816
$ws.Range("M34").Value = -1272.4
$ws.Range("N34").Value = -5194.1816

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 20385.375
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 20385.375
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 20385.375
$ws.Range("N41").Value = -21241.375

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1814.2222
$ws.Range("I132").Value = 1982.2667
$ws.Range("J132").Value = 1604.1666
$ws.Range("K132").Value = 5946.800099999999
$ws.Range("L132").Value = 4812.4998
$ws.Range("M132").Value = -3416.800099999999
$ws.Range("N132").Value = -9872.4998

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 121.296295
$ws.Range("I4").Value = 99
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 297
$ws.Range("L4").Value = 1200
$ws.Range("M4").Value = -185
$ws.Range("N4").Value = -1424

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 668056.8
$ws.Range("I121").Value = 50
$ws.Range("J121").Value = 703215.0600000001
$ws.Range("K121").Value = 150
$ws.Range("L121").Value = 2109645.18
$ws.Range("M121").Value = 1160
$ws.Range("N121").Value = -2112265.18

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 873.64
$ws.Range("I131").Value = 544.2857
$ws.Range("J131").Value = 927.2558
$ws.Range("K131").Value = 1632.8571
$ws.Range("L131").Value = 2781.7674
$ws.Range("M131").Value = 3407.1429
$ws.Range("N131").Value = -12861.7674

# GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8500
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 8500
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 8500
$ws.Range("N5").Value = -8724
$ws.Range("M5").ClearContents()

# GSM row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 34319.445
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 34319.445
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 34319.445
$ws.Range("N123").Value = -39219.445

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2688.1875
$ws.Range("I132").Value = 2501.3
$ws.Range("J132").Value = 2999.6667
$ws.Range("K132").Value = 7503.900000000001
$ws.Range("L132").Value = 8999.000100000001
$ws.Range("M132").Value = -4973.900000000001
$ws.Range("N132").Value = -14059.0001

# LTW row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 7268.8237
$ws.Range("I2").Value = 1500
$ws.Range("J2").Value = 9672.5
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 9672.5
$ws.Range("M2").Value = -1388
$ws.Range("N2").Value = -9896.5

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1289.2273
$ws.Range("I22").Value = 1041.3
$ws.Range("J22").Value = 1495.8334
$ws.Range("K22").Value = 1041.3
$ws.Range("L22").Value = 1495.8334
$ws.Range("M22").Value = -746.3
$ws.Range("N22").Value = -2085.8334

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1289.2273
$ws.Range("I27").Value = 1041.3
$ws.Range("J27").Value = 1495.8334
$ws.Range("K27").Value = 1041.3
$ws.Range("L27").Value = 1495.8334
$ws.Range("M27").Value = -934.3
$ws.Range("N27").Value = -1709.8334

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5747.7393
$ws.Range("I132").Value = 6036.8423
$ws.Range("J132").Value = 4374.5
$ws.Range("K132").Value = 18110.5269
$ws.Range("L132").Value = 13123.5
$ws.Range("M132").Value = -15580.5269
$ws.Range("N132").Value = -18183.5

# WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2000
$ws.Range("N2").Value = -2224

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3398.2144
$ws.Range("I132").Value = 4732
$ws.Range("J132").Value = 2397.875
$ws.Range("K132").Value = 14196
$ws.Range("L132").Value = 7193.625
$ws.Range("M132").Value = -11666
$ws.Range("N132").Value = -12253.625

Write-Host "Edit complete"
